$wb = $excel.ActiveWorkbook

# --- Add the new sheet "MaxRuntimeFirst" right after Sheet1 --------------
$sheet1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "MaxRuntimeFirst"

# --- Fill in the threshold / runtime results grid (A1:J4) ----------------
$data = @(
    @(5442282.407409668, 0.2, 0.3, 0.4, 0.5, 0.60000000000000009, 0.70000000000000007, 0.8, 0.9, $null),
    @(94973.148147583008, 101754362.7781982, 103439977.22257081, 104039521.85169069, 104702866.6664001, 102086890.184967, 103904323.1474915, 104725303.33363651, 102799408.7034851, 103921255.37055659),
    @(11745880.555557249, 104101626.1104065, 102955019.444519, 103506021.8514282, 105319586.8520508, 103442001.8518005, 105554180.9261536, 103402012.2221252, 103384954.073938, 104236394.4444702),
    @($null, 101889940.9261353, 104131162.7781616, 103218751.851355, 103204377.77811889, 101136545.55504151, 104416279.4444824, 103793333.70430911, 103735048.8889709, 101922700.0005554)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $val = $data[$r][$c]
        if ($null -ne $val) {
            $ws2.Cells.Item($r + 1, $c + 1).Value = $val
        }
    }
}

# --- Match the new sheet's page margins (1in top/bottom, 0.75in sides, 0.5in header/footer)
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# --- Make the new sheet the active / selected tab, with M6 selected ------
$ws2.Activate()
$null = $ws2.Range("M6").Select()
